# Updates the vendor config sheet from the "Jenne" (xlsx/admin-ftp) source
# settings to the new "Jenne - old" (csv/public-ftp) source settings, and
# turns the Source Ftp Url cell into a live hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lookup / Changes column headers for this vendor row were free-text
#     labels; they become the new source-file column names. These two
#     cells also lose their inherited banner styling (back to Normal). ---
$ws.Range("C2").Value = "ManPartNum"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "ListPr"
$ws.Range("D2").Style = "Normal"

# --- File extensions: xlsx source -> csv source ---
$ws.Range("E2").Value = ".csv"
$ws.Range("F2").Value = ".csv"

# --- Old file sheet name / postfix ---
$ws.Range("G2").Value = "Jenne - old"
$ws.Range("I2").Value = "JEN"

# --- Source FTP connection details move from the internal telquest FTP
#     to the vendor's own public FTP server. ---
$ws.Range("J2").Value = "ftp.jenne.com"
$ws.Range("K2").Value = "TE07004"
$ws.Range("L2").Value = "QQJjd2345"
$ws.Range("N2").Value = "/"
$ws.Range("O2").Value = "pricing"

# Turn the Source Ftp Url cell into a real hyperlink (adds the Hyperlink
# cell style + the relationship automatically).
$ws.Hyperlinks.Add($ws.Range("J2"), "http://ftp.jenne.com/")

# Leave the cursor where the editor left it.
$ws.Range("J15").Select()
